$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# Give the merged-header cells C1/D1 an outline border (top+bottom, then
# top+bottom+right for the closing cell) instead of the all-sides box
# they inherited from style "1".
$c1 = $ws1.Range("C1")
$c1.Borders(8).LineStyle = 1
$c1.Borders(9).LineStyle = 1

$d1 = $ws1.Range("D1")
$d1.Borders(8).LineStyle = 1
$d1.Borders(9).LineStyle = 1
$d1.Borders(10).LineStyle = 1

# Anonymize "fedcore" -> "approach" in the header row.
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

$c1b = $ws2.Range("C1")
$c1b.Borders(8).LineStyle = 1
$c1b.Borders(9).LineStyle = 1

$d1b = $ws2.Range("D1")
$d1b.Borders(8).LineStyle = 1
$d1b.Borders(9).LineStyle = 1
$d1b.Borders(10).LineStyle = 1

$f1b = $ws2.Range("F1")
$f1b.Borders(8).LineStyle = 1
$f1b.Borders(9).LineStyle = 1

$g1b = $ws2.Range("G1")
$g1b.Borders(8).LineStyle = 1
$g1b.Borders(9).LineStyle = 1
$g1b.Borders(10).LineStyle = 1

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell at G5.
$ws2.Range("G5").ClearContents()
